# Update "想去人数" (number of interested attendees) figures that changed
# between the previous data pull and the one generated at commit 456a3b4.
#
# Both the "展览" sheet and the aggregate "全部类型" sheet list the same
# exhibitions, so each event's updated count has to be written twice (the
# row numbers differ by one because "全部类型" has an extra row that isn't
# present in "展览").

$wb = $excel.ActiveWorkbook

# -- "展览" sheet --------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value  = 6760   # was 6755
$wsExpo.Range("F4").Value  = 427    # was 426
$wsExpo.Range("F15").Value = 1615   # was 1614
$wsExpo.Range("F21").Value = 2025   # was 2022
$wsExpo.Range("F22").Value = 141    # was 140

# -- "全部类型" sheet -----------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 6760    # was 6755
$wsAll.Range("F4").Value  = 427     # was 426
$wsAll.Range("F16").Value = 1615    # was 1614
$wsAll.Range("F22").Value = 2025    # was 2022
$wsAll.Range("F23").Value = 141     # was 140
